$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Row 5 (Georgia) ---
$ws.Range("D5").Value = 33508
$ws.Range("E5").Value = 1405
$ws.Range("F5").Value = 11857
$ws.Range("H5").Value = 35.39
$ws.Range("I5").Value = 49.61

# --- Update Row 6 (Michigan) ---
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value2 = "5/10/2020"
$ws.Range("D6").Value = 47138
$ws.Range("E6").Value = 4551
$ws.Range("F6").Value = 15084
$ws.Range("G6").Value = 1866

# --- Update Row 10 (Wisconsin -- Milwaukee) ---
$ws.Range("D10").Value = 3981
$ws.Range("E10").Value = 217
$ws.Range("F10").Value = 1544
$ws.Range("G10").Value = 106
$ws.Range("H10").Value = 38.78
$ws.Range("I10").Value = 48.85

# --- New Row 11 (San Diego) ---
$ws.Range("A11").Value = "San Diego"
$ws.Range("A11").Font.Bold = $true
$ws.Range("A11").HorizontalAlignment = -4108
$ws.Range("A11").VerticalAlignment = -4160
$ws.Range("A11").Borders.LineStyle = 1
$ws.Range("B11").Value = "California - San Diego"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value2 = "5/9/2020"
$ws.Range("D11").Value = 4926
$ws.Range("E11").Value = 175
$ws.Range("F11").Value = 167
$ws.Range("G11").Value = 4
$ws.Range("H11").Value = 3.39
$ws.Range("I11").Value = 2.29
$ws.Range("J11").Value = "Success!"

# --- New Row 12 (Florida) ---
$ws.Range("A12").Value = "Florida"
$ws.Range("A12").Font.Bold = $true
$ws.Range("A12").HorizontalAlignment = -4108
$ws.Range("A12").VerticalAlignment = -4160
$ws.Range("A12").Borders.LineStyle = 1
$ws.Range("B12").Value = "Florida"
$ws.Range("J12").Value = "An error occured. ... FileNotFoundError(2, 'No such file or directory')"
